# Applies the "cryptos list" data refresh described by the commit
# "Updated cryptos list on Thu Feb 29 05:07:48 UTC 2024 with GitHub Actions".
#
# The sheet stores every data cell (Coin/Link/Price/Volume) as literal text
# (inline strings in the XML), including "Price" values that look like plain
# numbers (e.g. 414.63) or clean percentages (e.g. 5.50%). Excel's normal
# parser would silently convert a bare numeric-looking assignment into an
# actual Number (dropping trailing zeros, e.g. "52.80" -> 52.8, or "1.00" ->
# 1), which would NOT match the source data. To keep those specific cells as
# genuine text - exactly like typing '414.63 into Excel - a leading apostrophe
# (the standard "treat as text" quote prefix) is used for them. Percentage
# cells already carry padding spaces/sign (e.g. "  +8.74%  ") so Excel leaves
# those as text on its own.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.272.57"
$ws.Range("E2").Value = "  +8.74%  "
$ws.Range("D3").Value = "3.437.66"
$ws.Range("E3").Value = "  +5.31%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'414.63"
$ws.Range("E5").Value = "  +4.24%  "
$ws.Range("D6").Value = "'123.18"
$ws.Range("E6").Value = "  +13.22%  "
$ws.Range("D7").Value = "3.432.40"
$ws.Range("E7").Value = "  +5.32%  "
$ws.Range("D8").Value = "'0.591"
$ws.Range("E8").Value = "  +1.69%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").Value = "'0.656"
$ws.Range("E10").Value = "  +5.48%  "
$ws.Range("D11").Value = "'0.129"
$ws.Range("E11").Value = "  +34.52%  "
$ws.Range("D12").Value = "'41.32"
$ws.Range("E12").Value = "  +4.67%  "
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").Value = "3.989.84"
$ws.Range("E14").Value = "  +5.51%  "
$ws.Range("D15").Value = "'8.51"
$ws.Range("E15").Value = "  +2.79%  "
$ws.Range("D16").Value = "'19.74"
$ws.Range("E16").Value = "  +3.81%  "
$ws.Range("D17").Value = "3.435.77"
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").Value = "62.340.20"
$ws.Range("E18").Value = "  +9.28%  "
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("D20").Value = "'10.84"
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").Value = "'0.0000132"
$ws.Range("E21").Value = "  +21.35%  "
$ws.Range("D22").Value = "'3.32"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").Value = "'83.79"
$ws.Range("E23").Value = "  +12.83%  "
$ws.Range("D24").Value = "'320.61"
$ws.Range("E24").Value = "  +8.93%  "
$ws.Range("D25").Value = "'12.96"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").Value = "'3.18"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "'31.14"
$ws.Range("E27").Value = "  +10.61%  "
$ws.Range("D28").Value = "'7.89"
$ws.Range("E28").Value = "  +6.65%  "
$ws.Range("D29").Value = "'7.83"
$ws.Range("E29").Value = "  -1.78%  "
$ws.Range("D30").Value = "'4.29"
$ws.Range("E30").Value = "  -2.05%  "
$ws.Range("E31").Value = "  +2.94%  "
$ws.Range("E32").Value = "  +4.13%  "
$ws.Range("E33").Value = "  +20.56%  "
# Rows 34/35 swapped position: InjectiveProtocol <-> Cosmos
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").Value = "'11.51"
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'42.04"
$ws.Range("E35").Value = "  +4.77%  "
$ws.Range("D37").Value = "'0.0484"
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("D38").Value = "'52.80"
$ws.Range("E38").Value = "  +2.79%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("E40").Value = "  +0.93%  "
$ws.Range("D41").Value = "'3.03"
$ws.Range("E41").Value = "  +0.32%  "
# Rows 42/43 swapped position: ARBITRUM <-> Stellar
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.125"
$ws.Range("E42").Value = "  +3.27%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'1.98"
$ws.Range("E43").Value = "  +5.76%  "
$ws.Range("D44").Value = "'134.45"
$ws.Range("E44").Value = "  -1.74%  "
$ws.Range("D45").Value = "'17.14"
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("D46").Value = "'0.283"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("D47").Value = "'3.89"
$ws.Range("E47").Value = "  -1.15%  "
# Rows 48/49 swapped position: WEMIXToken <-> EnergySwap
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'22.10"
$ws.Range("E48").Value = "  -1.16%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'2.21"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").Value = "2.201.07"
$ws.Range("E50").Value = "  +2.28%  "
$ws.Range("D51").Value = "3.785.83"
$ws.Range("E51").Value = "  +5.50%  "
